$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPart = 2
$xlByRows = 1

$used = $ws.UsedRange
$used.Replace("D64", "D69", $xlPart, $xlByRows, $true, $false, $false)
$used.Replace("D80", "D86", $xlPart, $xlByRows, $true, $false, $false)
$used.Replace("D51", "D55", $xlPart, $xlByRows, $true, $false, $false)
$used.Replace("S30", "S31", $xlPart, $xlByRows, $true, $false, $false)
